$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new log entry in row 9 (continuation of work on the background section)
$ws.Range("B9").Value = 'Some more work on "Background and Methodologies"'
$ws.Range("C9").Value = 0.5

# Move/leave the active selection where the user ended up after editing (D13)
$ws.Range("D13").Select()

$wb.Save()
